$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Unit Sales"
$ws.Range("A1:B1").Style = "Heading 1"
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("B1").Select() | Out-Null
